$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.721.47'
$ws.Range('D3').Value = '3.541.16'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('D7').Value = '3.542.04'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.47%  '
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('E12').Value = '  -2.93%  '
$ws.Range('D13').Value = '4.141.66'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000198'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('D16').Value = '3.542.11'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').Value = '65.683.14'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '418.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '3.679.89'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  -4.59%  '
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('E29').Value = '  -3.54%  '
$ws.Range('E30').Value = '  -4.84%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = '3.550.55'
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -8.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '174.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0816'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.50%  '
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.860'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('E45').Value = '  -7.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('E47').Value = '  -7.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.907'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.11%  '
